$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3917.349
$ws.Range("I64").Value = 3751.5151
$ws.Range("K64").Value = 3751.5151
$ws.Range("M64").Value = -3503.5151

$ws.Range("H67").Value = 3917.349
$ws.Range("I67").Value = 3751.5151
$ws.Range("K67").Value = 3751.5151
$ws.Range("M67").Value = -2893.5151

$ws.Range("H76").Value = 3256.9412
$ws.Range("I76").Value = 3112.606
$ws.Range("K76").Value = 3112.606
$ws.Range("M76").Value = -2797.606

$ws.Range("H79").Value = 3256.9412
$ws.Range("I79").Value = 3112.606
$ws.Range("K79").Value = 3112.606
$ws.Range("M79").Value = -2020.606

$ws.Range("H116").Value = 2432.1072
$ws.Range("J116").Value = 2342.8572
$ws.Range("L116").Value = 2342.8572
$ws.Range("N116").Value = -9226.8572

$ws.Range("H138").Value = 2863.69
$ws.Range("I138").Value = 1093.4615
$ws.Range("J138").Value = 3128.2068
$ws.Range("K138").Value = 3280.3845
$ws.Range("L138").Value = 9384.6204
$ws.Range("M138").Value = 1859.6155
$ws.Range("N138").Value = -19664.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2527
$ws.Range("I63").Value = 2568.6
$ws.Range("J63").Value = 2475
$ws.Range("K63").Value = 2568.6
$ws.Range("L63").Value = 2475
$ws.Range("M63").Value = -1882.6
$ws.Range("N63").Value = -3847

$ws.Range("H66").Value = 2527
$ws.Range("I66").Value = 2568.6
$ws.Range("J66").Value = 2475
$ws.Range("K66").Value = 12843
$ws.Range("L66").Value = 12375
$ws.Range("M66").Value = -9411
$ws.Range("N66").Value = -19239

$ws.Range("H74").Value = 1139.9048
$ws.Range("I74").Value = 971.9459000000001
$ws.Range("J74").Value = 2382.8
$ws.Range("K74").Value = 971.9459000000001
$ws.Range("L74").Value = 2382.8
$ws.Range("M74").Value = -97.94590000000005
$ws.Range("N74").Value = -4130.8

$ws.Range("H77").Value = 1139.9048
$ws.Range("I77").Value = 971.9459000000001
$ws.Range("J77").Value = 2382.8
$ws.Range("K77").Value = 4859.7295
$ws.Range("L77").Value = 11914
$ws.Range("M77").Value = -491.7295000000004
$ws.Range("N77").Value = -20650

$ws.Range("H132").Value = 1159.3469
$ws.Range("I132").Value = 592.4286
$ws.Range("J132").Value = 4560.857
$ws.Range("K132").Value = 1777.2858
$ws.Range("L132").Value = 13682.571
$ws.Range("M132").Value = 752.7142000000001
$ws.Range("N132").Value = -18742.571

$ws.Range("H137").Value = 43800
$ws.Range("J137").Value = 43800
$ws.Range("L137").Value = 43800
$ws.Range("N137").Value = -54000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12414.833
$ws.Range("I82").Value = 1322.25
$ws.Range("J82").Value = 34600
$ws.Range("K82").Value = 1322.25
$ws.Range("L82").Value = 34600
$ws.Range("M82").Value = -939.25
$ws.Range("N82").Value = -35366

$ws.Range("H85").Value = 12414.833
$ws.Range("I85").Value = 1322.25
$ws.Range("J85").Value = 34600
$ws.Range("K85").Value = 1322.25
$ws.Range("L85").Value = 34600
$ws.Range("M85").Value = 3.75
$ws.Range("N85").Value = -37252

$ws.Range("H105").Value = 2086.05
$ws.Range("I105").Value = 1779.2
$ws.Range("J105").Value = 3006.6
$ws.Range("K105").Value = 1779.2
$ws.Range("L105").Value = 3006.6
$ws.Range("M105").Value = -32.20000000000005
$ws.Range("N105").Value = -6500.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15220
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 16355.556
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 16355.556
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -17605.556

$ws.Range("H51").Value = 109280
$ws.Range("I51").Value = 502500
$ws.Range("J51").Value = 10975
$ws.Range("K51").Value = 502500
$ws.Range("L51").Value = 10975
$ws.Range("M51").Value = -501764
$ws.Range("N51").Value = -12447

$ws.Range("H58").Value = 1846.1628
$ws.Range("I58").Value = 1208.6207
$ws.Range("J58").Value = 3166.7856
$ws.Range("K58").Value = 1208.6207
$ws.Range("L58").Value = 3166.7856
$ws.Range("M58").Value = -1005.6207
$ws.Range("N58").Value = -3572.7856

$ws.Range("H59").Value = 12489.357
$ws.Range("J59").Value = 12489.357
$ws.Range("L59").Value = 12489.357
$ws.Range("N59").Value = -14779.357

$ws.Range("H61").Value = 109280
$ws.Range("I61").Value = 502500
$ws.Range("J61").Value = 10975
$ws.Range("K61").Value = 502500
$ws.Range("L61").Value = 10975
$ws.Range("M61").Value = -502152
$ws.Range("N61").Value = -11671

$ws.Range("H62").Value = 2715.3845
$ws.Range("I62").Value = 2516.6667
$ws.Range("J62").Value = 2885.7144
$ws.Range("K62").Value = 2516.6667
$ws.Range("L62").Value = 2885.7144
$ws.Range("M62").Value = -1892.6667
$ws.Range("N62").Value = -4133.7144

$ws.Range("H65").Value = 2715.3845
$ws.Range("I65").Value = 2516.6667
$ws.Range("J65").Value = 2885.7144
$ws.Range("K65").Value = 12583.3335
$ws.Range("L65").Value = 14428.572
$ws.Range("M65").Value = -9463.333500000001
$ws.Range("N65").Value = -20668.572

$ws.Range("H136").Value = 1846.1628
$ws.Range("I136").Value = 1208.6207
$ws.Range("J136").Value = 3166.7856
$ws.Range("K136").Value = 3625.8621
$ws.Range("L136").Value = 9500.356800000001
$ws.Range("M136").Value = -1075.8621
$ws.Range("N136").Value = -14600.3568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6877.2104
$ws.Range("I70").Value = 7215.706
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 7215.706
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -6945.706
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 6877.2104
$ws.Range("I73").Value = 7215.706
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 7215.706
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -6279.706
$ws.Range("N73").Value = -5872

$ws.Range("H80").Value = 2879.8
$ws.Range("I80").Value = 2667.3157
$ws.Range("J80").Value = 3035.077
$ws.Range("K80").Value = 2667.3157
$ws.Range("L80").Value = 3035.077
$ws.Range("M80").Value = -1669.3157
$ws.Range("N80").Value = -5031.077

$ws.Range("H83").Value = 2879.8
$ws.Range("I83").Value = 2667.3157
$ws.Range("J83").Value = 3035.077
$ws.Range("K83").Value = 13336.5785
$ws.Range("L83").Value = 15175.385
$ws.Range("M83").Value = -8344.5785
$ws.Range("N83").Value = -25159.385
